$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (row 5 and row 6)
$ws.Range("D5").Value = 40

$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 3

# Update the view: zoom to 130% and move the selection to D5
$excel.ActiveWindow.Zoom = 130
$ws.Range("D5").Select()
